# Updated cryptos list on Thu Jan 11 19:30:49 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.401.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.605.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.81%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.51%  "

# Row 7
$ws.Range("E7").Value = "  +4.27%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.23%  "

# Row 9
$ws.Range("E9").Value = "  +12.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0847"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.48%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.003.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.55%  "

# Row 15
$ws.Range("E15").Value = "  +1.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.613.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.05%  "

# Row 17
$ws.Range("E17").Value = "  +7.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.507.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "

# Row 20
$ws.Range("E20").Value = "  +7.12%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.90%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.65%  "

# Row 25
$ws.Range("E25").Value = "  +7.71%  "

# Row 26
$ws.Range("E26").Value = "  +9.32%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +32.28%  "

# Row 28
$ws.Range("E28").Value = "  -0.32%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.02"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.65%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.60%  "

# Row 32
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "39.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.00%  "

# Row 34
$ws.Range("E34").Value = "  -6.82%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0843"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.27%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.42%  "

# Row 37
$ws.Range("E37").Value = "  +7.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.94%  "

# Row 39
$ws.Range("E39").Value = "  +7.63%  "

# Row 40
$ws.Range("E40").Value = "  +5.51%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +38.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.11%  "

# Row 43
$ws.Range("E43").Value = "  +9.24%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.71%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.45%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.131.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.77%  "

# Row 47
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "93.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.25%  "

# Row 50
$ws.Range("E50").Value = "  -1.47%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.74%  "
